# Applies the "#12 Finish with request" edit to Protocol.docx:
#  1. Fill in the last real table row (previously all "?") describing the
#     "Server Rejected Request" message.
#  2. Move the stray "_GoBack" bookmark from the first paragraph to the
#     freshly-filled "Side" cell (Word re-drops _GoBack at the last edit
#     point whenever the document is saved).
#  3. Re-flow two bullet paragraphs so the spell-checker's
#     proofErr (spellStart/spellEnd) markers bracket the words that were
#     flagged ("т.ч" and "одмена"), matching a real Word re-save after
#     the text was touched.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark sitting after "... байта, значение"
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Fill the "Server Rejected Request" row (the row whose four cells
#    still hold the "?" placeholder) and re-drop the bookmark in the
#    last ("Side") cell, right before its new run.
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)

for ($i = 1; $i -le $table.Rows.Count; $i++) {
    if ($table.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13) -eq "?") {
        $targetRow = $i
    }
}

$table.Cell($targetRow, 1).Range.Text = "Server Rejected Request"
$table.Cell($targetRow, 2).Range.Text = "0xB"
$table.Cell($targetRow, 3).Range.Text = "reason"
$table.Cell($targetRow, 4).Range.Text = "Server"

$sideCellStart = $table.Cell($targetRow, 4).Range.Start
$collapsed = $d.Range($sideCellStart, $sideCellStart)
[void]$d.Bookmarks.Add("_GoBack", $collapsed)

# ---------------------------------------------------------------------
# 3) Re-split the two bullet paragraphs so proofErr markers wrap the
#    words the spell checker flagged, without changing the visible text.
# ---------------------------------------------------------------------
$wdParagraph = 4

$find1 = $d.Content
$find1.Find.ClearFormatting()
$found1 = $find1.Find.Execute("Обработка ответов сервера (в т.ч. чужие сообщения)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    [void]$find1.Expand($wdParagraph)
    $xml1 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:r><w:t xml:space='preserve'>Обработка ответов сервера (в </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r><w:t>т.ч</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r><w:t>. чужие сообщения)</w:t></w:r>" + `
        "</w:p>"
    [void]$find1.InsertXML($xml1)
}

$find2 = $d.Content
$find2.Find.ClearFormatting()
$found2 = $find2.Find.Execute("Инструменты одмена (кик, что-либо ещё)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    [void]$find2.Expand($wdParagraph)
    $xml2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:r><w:t xml:space='preserve'>Инструменты </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r><w:t>одмена</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r><w:t xml:space='preserve'> (кик, что-либо ещё)</w:t></w:r>" + `
        "</w:p>"
    [void]$find2.InsertXML($xml2)
}
